$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Section: Phase III header ----
$ws.Cells.Item(97,1).Value = "Phase III"
$ws.Cells.Item(97,1).Font.Bold = $true

# ---- Moni block ----
$ws.Cells.Item(98,1).Value = "Moni"
$ws.Cells.Item(99,1).Value = "Widrawl + returns"
$ws.Cells.Item(99,3).Value = 20630
$ws.Cells.Item(100,1).Value = "Paint"
$ws.Cells.Item(100,3).Value = 7900
$ws.Cells.Item(101,1).Value = "Wood material"
$ws.Cells.Item(101,3).Value = 3375
$ws.Cells.Item(102,1).Value = "Sultaan"
$ws.Cells.Item(102,3).Value = 6500
$ws.Cells.Item(103,1).Value = "Keelen"
$ws.Cells.Item(103,3).Value = 220
$ws.Cells.Item(104,1).Value = "Lock (sultan)"
$ws.Cells.Item(104,3).Value = 165
$ws.Cells.Item(105,1).Value = "Material"
$ws.Cells.Item(105,3).Value = 400
$ws.Cells.Item(106,1).Value = "Paint operations"
$ws.Cells.Item(106,3).Value = 1150
$ws.Cells.Item(107,3).Value = 500
$ws.Cells.Item(108,1).Value = "Total cost"
$ws.Cells.Item(108,3).Formula = "=SUM(C100:C107)"
$ws.Cells.Item(109,1).Value = "Remaining amount"
$ws.Cells.Item(109,3).Formula = "=C99-C108"

# ---- Sultan block ----
$ws.Cells.Item(113,1).Value = "Sultan"
$ws.Cells.Item(114,1).Value = "Paint"
$ws.Cells.Item(114,3).Value = 1900
$ws.Cells.Item(115,1).Value = "labour"
$ws.Cells.Item(115,3).Value = 700
$ws.Cells.Item(116,1).Value = "paint material"
$ws.Cells.Item(116,3).Value = 100
$ws.Cells.Item(117,1).Value = "seeri"
$ws.Cells.Item(117,3).Value = 1000
$ws.Cells.Item(118,1).Value = "flex"
$ws.Cells.Item(118,3).Value = 1300
$ws.Cells.Item(119,1).Value = "Donkey cart"
$ws.Cells.Item(119,3).Value = 150
$ws.Cells.Item(120,1).Value = "Carpenter"
$ws.Cells.Item(120,3).Value = 800
$ws.Cells.Item(121,1).Value = "Zink"
$ws.Cells.Item(121,3).Value = 900
$ws.Cells.Item(122,1).Value = "grey teen paint"
$ws.Cells.Item(122,3).Value = 700
$ws.Cells.Item(123,1).Value = "labour "
$ws.Cells.Item(123,3).Value = 5200
$ws.Cells.Item(126,1).Value = "Total"
$ws.Cells.Item(126,3).Formula = "=SUM(C114:C125)"

$ws.Cells.Item(107,1).Value = "inouguration+Van Rent"

$ws.Range("C108").Select()
